# Sprint 4 Backlog - Burndown: record Janera's actual time spent on the
# "Add planned meal ingredients to shopping list" task (row 6) and the
# associated Week-1 "Amount Remaining" hours for that row plus two other
# in-flight tasks (rows 16 and 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "Add planned meal ingredients to shopping list" /
# "Create functionality to add needed ingredients from planned meals to
# shpping list (desktop)" — log 1.5 hrs actual time, completed by Janera,
# and reflect that in the Week 1 / Week 2 remaining-time columns.
$ws.Range("E6").Value = 1.5
$ws.Range("F6").Value = "Janera"
$ws.Range("H6").Value = 1.5
$ws.Range("I6").Value = 0

# Row 16 and Row 22: additional Week 1 time logged this sprint.
$ws.Range("H16").Value = 0.5
$ws.Range("H22").Value = 2

# Update the active selection to match where the author left off editing.
$ws.Range("G15").Select() | Out-Null
